# individual leads import done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SERVICE_NAME column (column D) is no longer part of the individual
# leads import -- delete it entirely so every column to its right shifts
# left by one.
$ws.Columns("D").Delete()

# Fix the mis-spelled header: DATA_OF_BIRTH -> DATE_OF_BIRTH
$ws.Range("F1").Value = "DATE_OF_BIRTH"

# Update the active cell / selection to match where editing left off.
$ws.Range("F6").Select()
